$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad / "Changed") date value from 45203 to 45205
# for every data row (rows 2 through 171).
$ws.Range("C2:C171").Value = 45205
